$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values (B2:E2)
$ws.Range("B2").Value = 11.73003337338305
$ws.Range("C2").Value = 39.219110102977254
$ws.Range("D2").Value = 51.827360856362169
$ws.Range("E2").Value = 44.6426690482119

# Row 3 values (B3, C3 cleared, D3 new, E3)
$ws.Range("B3").Value = 29.667566400003114
$null = $ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 65.142064585136268
$ws.Range("E3").Value = 43.378781294180513

# Reflect the edited range as the current selection
$null = $ws.Range("B1:E3").Select()
